$d = $word.ActiveDocument

# --- 1. Merge the "Fridlysta arter" intro paragraph with the bullet item
#        that followed it: the species name becomes lower-case and the
#        whole thing turns into one normal-style sentence ending in a
#        period; the separate bulleted paragraph disappears. ---
$introText = "Följande fridlysta arter har sina livsmiljöer och växtplatser i den avverkningsanmälda skogen: fläcknycklar (§8)."

$found = $false
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Följande fridlysta arter har sina livsmiljöer")) {
        $introPara = $p
        $found = $true
        break
    }
}

if ($found) {
    $bulletPara = $introPara.Next()
    $introPara.Range.Text = $introText
    $bulletPara.Range.Delete()
}

# --- 2. Bump the date shown in the document's first-page header. ---
foreach ($story in $d.StoryRanges) {
    $rng = $story
    while ($rng -ne $null) {
        $rng.Find.Execute("2023-10-22", $true, $false, $false, $false, $false,
                           $true, 1, $false, "2023-10-25", 2) | Out-Null
        $rng = $rng.NextStoryRange
    }
}
